$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.714.87"
$ws.Range("E2").Value = "  +0.26%  "

$ws.Range("D3").Value = "1.601.15"
$ws.Range("E3").Value = "  +0.28%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.12%  "

$ws.Range("E6").Value = "  -0.39%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.01"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("E9").Value = "  +0.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.77%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0846"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.00%  "

$ws.Range("D12").Value = "1.826.16"
$ws.Range("E12").Value = "  +0.30%  "

$ws.Range("D13").Value = "1.601.96"
$ws.Range("E13").Value = "  -0.03%  "

$ws.Range("E14").Value = "  +0.50%  "

$ws.Range("E15").Value = "  -0.25%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.03%  "

$ws.Range("E17").Value = "  +0.36%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "210.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.58%  "

$ws.Range("E21").Value = "  -0.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.68"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.88%  "

$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("E26").Value = "  -0.14%  "

$ws.Range("E27").Value = "  -0.92%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.64%  "

$ws.Range("E29").Value = "  -0.78%  "

$ws.Range("E30").Value = "  +0.01%  "

$ws.Range("E31").Value = "  +0.84%  "

$ws.Range("D33").Value = "1.291.03"
$ws.Range("E33").Value = "  +0.71%  "

$ws.Range("E34").Value = "  +0.56%  "

$ws.Range("E35").Value = "  +0.54%  "

$ws.Range("E36").Value = "  -2.68%  "

$ws.Range("E37").Value = "  +10.72%  "

$ws.Range("E38").Value = "  -0.38%  "

$ws.Range("E39").Value = "  -0.42%  "

$ws.Range("E40").Value = "  -1.95%  "

$ws.Range("E41").Value = "  +0.18%  "

$ws.Range("E42").Value = "  +0.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "62.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.70%  "

$ws.Range("D44").Value = "1.737.71"
$ws.Range("E44").Value = "  +0.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.07%  "

$ws.Range("E46").Value = "  -1.53%  "

$ws.Range("E47").Value = "  -0.24%  "

$ws.Range("E48").Value = "  +1.50%  "

$ws.Range("E49").Value = "  +0.18%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.41%  "

$ws.Range("E51").Value = "  +0.90%  "
